# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (fund holdings for the quarter) right
# after "总计" and before "2022-Q1", and adds a matching summary row to the
# "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row 2 for 2022-Q3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Pull matching formatting (borders/font/alignment) for the new index
# cell from the row below it, then clear the incidental format that the
# row-insert auto-applied to the data cells.
$total.Range("A4").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.36

# The A column is a plain row-position index (0,1,2,...); renumber the
# two rows pushed down by the insert so it stays sequential.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 2. New "2022-Q3" sheet with the quarter's fund-holding detail.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Fetch the "2022-Q1" sheet reference AFTER the Add() call (adding a
# sheet reshuffles the collection, so a handle grabbed beforehand goes
# stale for clipboard ops) and copy its formatting (header style +
# index-column style) so the new sheet matches the workbook's look.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Range("B1:H1").Copy()
$q3.Range("B1").PasteSpecial(-4122)
$q1.Range("A2:H2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$q3.Range("A3").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "013991"
$q3.Range("B2").Style = "Normal"
$q3.Range("C2").Value = "中欧港股通精选一年持有混合A"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "6.69"
$q3.Range("D2").Style = "Normal"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "93.38"
$q3.Range("E2").Style = "Normal"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.21"
$q3.Range("F2").Style = "Normal"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.2147"
$q3.Range("G2").Style = "Normal"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "013992"
$q3.Range("B3").Style = "Normal"
$q3.Range("C3").Value = "中欧港股通精选一年持有混合C"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "4.68"
$q3.Range("D3").Style = "Normal"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "93.38"
$q3.Range("E3").Style = "Normal"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "3.21"
$q3.Range("F3").Style = "Normal"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.1502"
$q3.Range("G3").Style = "Normal"
$q3.Range("H3").Value = 10
